$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "SCRIPT/P01P04A/um1306.ssb" entry, same formatting as rows 4-7 (bordered style) ---
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A8:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = "SCRIPT/P01P04A/um1306.ssb"
$ws.Range("B8").Value = 279
$ws.Range("C8").Value = " I wonder what you get when you\nwin big?"
$ws.Range("D8").Value = " Интересно, что можно получить\nпри крупном выигрыше?"
$ws.Range("E8").Value = " Éîóåñåòîï, œóï íïçîï ðïìôœéóû\nðñé ëñôðîïí âúéãñúšå?"
$ws.Rows("8:8").RowHeight = 43.2

# --- Row 9: "SCRIPT/T01P01A/um1408.ssb" entry, same formatting as row 2 (unbordered style) ---
# (shared-string allocation order in the source file puts column C before column A here)
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Value = " For sure! Happy-happy! ♪[K]\nTee-hee! ♪"
$ws.Range("A9").Value = "SCRIPT/T01P01A/um1408.ssb"
$ws.Range("B9").Value = 260
$ws.Range("D9").Value = " В самом деле! Рады-рады! ♪[K]\nХи-хии! ♪"
$ws.Range("E9").Value = " Â òàíïí äåìå! Ñàäú-ñàäú! ♪[K]\nÖé-öéé! ♪"
$ws.Rows("9:9").RowHeight = 45.6

# --- Row 10: only column A filled ("SCRIPT/T01P01A/um1505.ssb"), same formatting as A9 (unbordered style) ---
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = "SCRIPT/T01P01A/um1505.ssb"
$ws.Rows("10:10").RowHeight = 43.2

# --- View state: scroll down and select C8, matching where the new rows were added ---
$ws.Range("C8").Select() | Out-Null
